$wb = $excel.ActiveWorkbook

# ALC!row 41
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 855.06665
$ws.Range("I41").Value = 1093.2727
$ws.Range("J41").Value = 200
$ws.Range("K41").Value = 1093.2727
$ws.Range("L41").Value = 200
$ws.Range("M41").Value = -653.2727
$ws.Range("N41").Value = -1080

# ALC!row 55
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 202.90244
$ws.Range("I55").Value = 201.92105
$ws.Range("J55").Value = 215.33333
$ws.Range("K55").Value = 201.92105
$ws.Range("L55").Value = 215.33333
$ws.Range("M55").Value = 12.07894999999999
$ws.Range("N55").Value = -643.3333299999999

# ALC!row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3531.9167
$ws.Range("I62").Value = 3113.3125
$ws.Range("J62").Value = 4369.125
$ws.Range("K62").Value = 3113.3125
$ws.Range("L62").Value = 4369.125
$ws.Range("M62").Value = -2489.3125
$ws.Range("N62").Value = -5617.125

# ALC!row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3531.9167
$ws.Range("I65").Value = 3113.3125
$ws.Range("J65").Value = 4369.125
$ws.Range("K65").Value = 15566.5625
$ws.Range("L65").Value = 21845.625
$ws.Range("M65").Value = -12446.5625
$ws.Range("N65").Value = -28085.625

# ALC!row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 63004.35
$ws.Range("I113").Value = 113452.22
$ws.Range("J113").Value = 6250.5
$ws.Range("K113").Value = 113452.22
$ws.Range("L113").Value = 6250.5
$ws.Range("M113").Value = -110198.22
$ws.Range("N113").Value = -12758.5

# ALC!row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3248.1177
$ws.Range("I132").Value = 2958.5715
$ws.Range("J132").Value = 4599.3335
$ws.Range("K132").Value = 8875.7145
$ws.Range("L132").Value = 13798.0005
$ws.Range("M132").Value = -6345.7145
$ws.Range("N132").Value = -18858.0005

# ALC!row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 91724.82000000001
$ws.Range("I135").Value = 717.25
$ws.Range("K135").Value = 6455.25
$ws.Range("M135").Value = -3920.25

# ALC!row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1405.1305
$ws.Range("I137").Value = 1215.3334
$ws.Range("J137").Value = 2088.4
$ws.Range("K137").Value = 3646.0002
$ws.Range("L137").Value = 6265.200000000001
$ws.Range("M137").Value = -1096.0002
$ws.Range("N137").Value = -11365.2

# ARM!row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 58.11111
$ws.Range("J5").Value = 100
$ws.Range("L5").Value = 100
$ws.Range("N5").Value = -324

# ARM!row 46
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 6204.6
$ws.Range("J46").Value = 6628.6665
$ws.Range("L46").Value = 6628.6665
$ws.Range("N46").Value = -7266.6665

# ARM!row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1583
$ws.Range("I74").Value = 1375.4736
$ws.Range("K74").Value = 1375.4736
$ws.Range("M74").Value = -501.4736

# ARM!row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1583
$ws.Range("I77").Value = 1375.4736
$ws.Range("K77").Value = 6877.368
$ws.Range("M77").Value = -2509.368

# ARM!row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2864.5217
$ws.Range("I132").Value = 1117.75
$ws.Range("K132").Value = 3353.25
$ws.Range("M132").Value = -823.25

# BSM!row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 58.11111
$ws.Range("J4").Value = 100
$ws.Range("L4").Value = 100
$ws.Range("N4").Value = -330

# BSM!row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 92294.73
$ws.Range("I107").Value = 144149.14
$ws.Range("J107").Value = 1549.5
$ws.Range("K107").Value = 144149.14
$ws.Range("L107").Value = 1549.5
$ws.Range("M107").Value = -142229.14
$ws.Range("N107").Value = -5389.5

# CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 41056.94
$ws.Range("I31").Value = 39746.85
$ws.Range("J31").Value = 45478.5
$ws.Range("K31").Value = 39746.85
$ws.Range("L31").Value = 45478.5
$ws.Range("M31").Value = -39451.85
$ws.Range("N31").Value = -46068.5

# CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 41056.94
$ws.Range("I34").Value = 39746.85
$ws.Range("J34").Value = 45478.5
$ws.Range("K34").Value = 39746.85
$ws.Range("L34").Value = 45478.5
$ws.Range("M34").Value = -39544.85
$ws.Range("N34").Value = -45882.5

# CRP!row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1818.8649
$ws.Range("I58").Value = 1784.9584
$ws.Range("K58").Value = 1784.9584
$ws.Range("M58").Value = -1581.9584

# CRP!row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1818.8649
$ws.Range("I136").Value = 1784.9584
$ws.Range("K136").Value = 5354.8752
$ws.Range("M136").Value = -2804.8752

# CUL!row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 406.54166
$ws.Range("J5").Value = 294.5
$ws.Range("L5").Value = 883.5
$ws.Range("N5").Value = -1107.5

# CUL!row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3624.5
$ws.Range("I70").Value = 3332.6667
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 9998.000100000001
$ws.Range("L70").Value = 13500
$ws.Range("M70").Value = -9683.000100000001
$ws.Range("N70").Value = -14130

# CUL!row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 3624.5
$ws.Range("I73").Value = 3332.6667
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 9998.000100000001
$ws.Range("L73").Value = 13500
$ws.Range("M73").Value = -8906.000100000001
$ws.Range("N73").Value = -15684

# CUL!row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1288.6666
$ws.Range("J113").Value = 1321.2941
$ws.Range("L113").Value = 3963.8823
$ws.Range("N113").Value = -8303.882300000001

# CUL!row 116
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 16120.571
$ws.Range("I116").Value = 18474
$ws.Range("K116").Value = 55422
$ws.Range("M116").Value = -51980

# CUL!row 118
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 3420
$ws.Range("I118").Value = 2433.3333
$ws.Range("K118").Value = 7299.999899999999
$ws.Range("M118").Value = -6056.999899999999

# CUL!row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 406.54166
$ws.Range("J135").Value = 294.5
$ws.Range("L135").Value = 2650.5
$ws.Range("N135").Value = -7720.5

# GSM!row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2350.389
$ws.Range("I126").Value = 2350.389
$ws.Range("K126").Value = 7051.167
$ws.Range("M126").Value = -4581.167

# GSM!row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3422.111
$ws.Range("I132").Value = 2770.3865
$ws.Range("J132").Value = 6289.7
$ws.Range("K132").Value = 8311.1595
$ws.Range("L132").Value = 18869.1
$ws.Range("M132").Value = -5781.1595
$ws.Range("N132").Value = -23929.1

# LTW!row 80
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# LTW!row 81
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# LTW!row 83
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# LTW!row 84
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# LTW!row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4254.769
$ws.Range("I136").Value = 3954.4443
$ws.Range("J136").Value = 4930.5
$ws.Range("K136").Value = 11863.3329
$ws.Range("L136").Value = 14791.5
$ws.Range("M136").Value = -9313.332900000001
$ws.Range("N136").Value = -19891.5

# WVR!row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 27779236
$ws.Range("I107").Value = 1895.7273
$ws.Range("J107").Value = 71429340
$ws.Range("K107").Value = 5687.1819
$ws.Range("L107").Value = 214288020
$ws.Range("M107").Value = -3767.1819
$ws.Range("N107").Value = -214291860

# WVR!row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4435.6855
$ws.Range("I132").Value = 5209.654
$ws.Range("J132").Value = 2199.7778
$ws.Range("K132").Value = 15628.962
$ws.Range("L132").Value = 6599.3334
$ws.Range("M132").Value = -13098.962
$ws.Range("N132").Value = -11659.3334
